# UOMModuleTest.xlsx: rename the "AutomationPR" UOM test data to "AutomationPQR"
# on the tc_UOM_3_4_5_7_9_10 sheet (B2 holds the UOM name, F2 holds the
# "...Updated" success message that is built from it).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tc_UOM_3_4_5_7_9_10")

$ws.Range("B2").Value = "AutomationPQR"
$ws.Range("F2").Value = "AutomationPQR Updated"
